# Shenango.pptx checkpoint edit:
#   1. Bump the cached "last saved" date field (datetimeFigureOut) on the
#      slide master and every slide layout from 2023. 12. 29. to 2023. 12. 30.
#   2. Fix the "What can be futher?" typo -> "What can be further?" on the
#      four closing slides that share that textbox.

$p = $ppt.ActivePresentation

$oldDate = "2023. 12. 29."
$newDate = "2023. 12. 30."

# --- 1a. Slide master date placeholder -------------------------------------
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 1b. Every slide layout's date placeholder ------------------------------
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Fix the "futher" -> "further" typo on every slide that has it ------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($i = 1; $i -le $sl.Shapes.Count; $i++) {
        $shp = $sl.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "What can be futher?") {
                $shp.TextFrame.TextRange.Text = "What can be further?"
            }
        }
    }
}
